# Add a new board sheet "Sheet1 (4)" to the ReversiKata workbook.
#
# The new sheet represents the next move in the game: it starts as a
# duplicate of the current last board sheet ("Sheet1 (3)") - same layout,
# styles, helper formulas (K:Q) and print setup - and then the move is
# applied: the two white discs at F7/G7 and the white disc at G8 are
# flipped to black, reflecting a black move that captured them.
#
# The new sheet is inserted right after "Sheet1 (3)" and becomes the
# active tab; "Sheet1 (3)" stops being the active tab.

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("Sheet1 (3)")

# Duplicate "Sheet1 (3)" and place the copy immediately after it.
$src.Copy($null, $src)
$newSheet = $wb.Worksheets.Item("Sheet1 (3) (2)")
$newSheet.Name = "Sheet1 (4)"

# Apply the move: flip the captured white ("W") discs to black ("B").
$newSheet.Range("F7").Value = "B"
$newSheet.Range("G7").Value = "B"
$newSheet.Range("G8").Value = "B"

# Restore/refresh the cursor position on the old sheet, then make the new
# sheet the active tab with the default selection.
$src.Activate()
$src.Range("M25").Select() | Out-Null

$newSheet.Activate()
$newSheet.Range("A1").Select() | Out-Null
